# This workbook lists daily price observations for "Achicoria" (rows 3-21).
# The commit re-shuffles the per-row data (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Origen, Precio $/Kg) across the
# existing rows 3-21, while row 2 and all the "constant" columns
# (A,B,C,E,F,G,H,I,N,Q,R) stay untouched.
#
# Mapping: new row R gets the data that used to live in row Src(R).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    3  = 10
    4  = 18
    5  = 3
    6  = 15
    7  = 12
    8  = 13
    9  = 5
    10 = 8
    11 = 20
    12 = 11
    13 = 19
    14 = 9
    15 = 6
    16 = 21
    17 = 7
    18 = 14
    19 = 17
    20 = 4
    21 = 16
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the current values of the columns that move, for every source row,
# before any writes happen (since several rows both give and receive data).
$snapshot = @{}
foreach ($r in 3..21) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# Now write each destination row using the snapshot of its mapped source row.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
